$d = $word.ActiveDocument

# 1) Merge the three runs "As my panic subsides, I start to notice a tinge of " +
#    "curiosity" + " in the back of my mind." into a single run/sentence.
$find1 = $d.Content.Find
$find1.Execute(
    "As my panic subsides, I start to notice a tinge of curiosity in the back of my mind.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As my panic subsides, I start to notice a tinge of curiosity in the back of my mind.",
    2)

# 2) Swap the stage direction from "neutral expressionless" to "neutral thinking".
$find2 = $d.Content.Find
$find2.Execute(
    "Mara (neutral expressionless): I wouldn" + [char]8217 + "t exactly say that. I left a little bit of bedhead for you.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Mara (neutral thinking): I wouldn" + [char]8217 + "t exactly say that. I left a little bit of bedhead for you.",
    2)
